# Updated data to reflect new requirement separation
#
# Inserts three new columns (Corequisites, Concurrent, Recommended) between
# the existing "Prerequisites" (C) and "Terms Typically Offered" (old D, now
# G) columns, fills them with "NA" for every course row, and tidies up a
# couple of Prerequisites descriptions that no longer mention a co-requisite
# explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D; this pushes the existing
# "Terms Typically Offered" column from D to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row values.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill the new columns with "NA" for every data row (2-31).
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# Update a few Prerequisites cells whose wording changed now that
# co-requisites are tracked separately. The course catalog text uses a
# non-breaking space between course prefixes and numbers (e.g. "CD 102").
$nbsp = [char]0x00A0
$ws.Range("C10").Value = "CD" + $nbsp + "102 and CD/PSY" + $nbsp + "256. Any 300-400 level CD course."
$ws.Range("C19").Value = "CD" + $nbsp + "329; and two of the CD" + $nbsp + "304, 305, or CD" + $nbsp + "306."
$ws.Range("C21").Value = "Two of the CD" + $nbsp + "304, CD" + $nbsp + "305, CD" + $nbsp + "306; junior standing or consent of instructor."
